$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels in row 1: "_old" -> "_FV2210", "_new" -> "_FV2304"
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the used range into a table (ListObject) named Table1
$rng = $ws.Range("A1:U72")
$listObject = $ws.ListObjects.Add(1, $rng, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# Freeze panes below the header row (row 2)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
